$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (date bumped from 04-05 to 04-06)
$ws.Name = "Through 2022-04-06"

# Update the "April (through 04-05)" label to "April (through 04-06)"
$ws.Range("A5").Value = "April (through 04-06)"

# Update April row (row 5) values
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 9
$ws.Range("F5").Value = 8
$ws.Range("G5").Value = 12
$ws.Range("H5").Value = 15
$ws.Range("I5").Value = 17

# Update Total row (row 6) values
$ws.Range("C6").Value = 133
$ws.Range("D6").Value = 198
$ws.Range("F6").Value = 118
$ws.Range("G6").Value = 210
$ws.Range("H6").Value = 438
$ws.Range("I6").Value = 451
